$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: set B column text values as formulas that evaluate to the literal
# string (including the leading apostrophe character), so the value is not
# auto-corrected into Excel's "quote prefix" text marker.
$ws.Cells.Item(2,2).Formula = "=""'Bacteroides_cellulosilyticus_DSM_14838.mat'"""
$ws.Cells.Item(3,2).Formula = "=""'Bacteroides_coprocola_M16_DSM_17136.mat'"""
$ws.Cells.Item(4,2).Formula = "=""'Bacteroides_coprophilus_DSM_18228.mat'"""
$ws.Cells.Item(5,2).Formula = "=""'Bacteroides_fluxus_YIT_12057.mat'"""
$ws.Cells.Item(6,2).Formula = "=""'Bacteroides_oleiciplenus_YIT_12058.mat'"""
$ws.Cells.Item(7,2).Formula = "=""'Bacteroides_ovatus_ATCC_8483.mat'"""
$ws.Cells.Item(8,2).Formula = "=""'Bacteroides_plebeius_M12_DSM_17135.mat'"""
$ws.Cells.Item(9,2).Formula = "=""'Bacteroides_salyersiae_WAL_10018.mat'"""
$ws.Cells.Item(10,2).Formula = "=""'Bacteroides_stercoris_ATCC_43183.mat'"""
$ws.Cells.Item(11,2).Formula = "=""'Bacteroides_thetaiotaomicron_VPI_5482.mat'"""
$ws.Cells.Item(12,2).Formula = "=""'Bacteroides_uniformis_ATCC_8492.mat'"""
$ws.Cells.Item(13,2).Formula = "=""'Bacteroides_vulgatus_ATCC_8482.mat'"""
$ws.Cells.Item(14,2).Formula = "=""'Bifidobacterium_animalis_lactis_AD011.mat'"""
$ws.Cells.Item(15,2).Formula = "=""'Enterococcus_faecalis_OG1RF_ATCC_47077.mat'"""
$ws.Cells.Item(16,2).Formula = "=""'Flavonifractor_plautii_ATCC_29863.mat'"""
$ws.Cells.Item(17,2).Formula = "=""'Lactobacillus_plantarum_JDM1.mat'"""
$ws.Cells.Item(18,2).Formula = "=""'Odoribacter_laneus_YIT_12061.mat'"""
$ws.Cells.Item(19,2).Formula = "=""'Parabacteroides_distasonis_ATCC_8503.mat'"""
$ws.Cells.Item(20,2).Formula = "=""'Parabacteroides_johnsonii_DSM_18315.mat'"""

# Step 2: convert the formulas above into static values via copy / paste-special
# (xlPasteValues = -4163), matching how the apostrophe-led text would be "baked in"
$bRange = $ws.Range("B2:B20")
$bRange.Copy()
$bRange.PasteSpecial(-4163)

# Step 3: update the C column numeric values
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(8,3).Value = 0
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(11,3).Value = 0.016
$ws.Cells.Item(12,3).Value = 0
$ws.Cells.Item(13,3).Value = 0.011
$ws.Cells.Item(14,3).Value = 0
$ws.Cells.Item(15,3).Value = 0
$ws.Cells.Item(16,3).Value = 0
$ws.Cells.Item(17,3).Value = 0.014
$ws.Cells.Item(18,3).Value = 0.144
$ws.Cells.Item(19,3).Value = 0
$ws.Cells.Item(20,3).Value = 0.8129999999999999

$excel.CutCopyMode = 0
